$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("D:D").Insert()

# Copy formatting from column E (the old D, now shifted) to new column D
$ws.Columns("E:E").Copy()
$ws.Columns("D:D").PasteSpecial(-4122)  # xlPasteFormats = -4122

Write-Host ("D7 NumberFormat: " + $ws.Range("D7").NumberFormat)
Write-Host ("D8 NumberFormat: " + $ws.Range("D8").NumberFormat)
Write-Host ("D7 Font Bold: " + $ws.Range("D7").Font.Bold)
Write-Host ("D8 Font Bold: " + $ws.Range("D8").Font.Bold)
